$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update task statuses (RRHH liquidacion de SAC: se quita el calculo de
# antiguedad -> los items correspondientes pasan a "terminado"/"en proceso").
$ws.Range("B101").Value = "en proceso"
$ws.Range("B115").Value = "terminado"
$ws.Range("B117").Value = "terminado"
$ws.Range("B131").Value = "terminado"
$ws.Range("B136").Value = "terminado"

# These rows no longer match the active AutoFilter criteria ("no comenzado"),
# so they become hidden just like the rest of the filtered-out rows.
$ws.Rows.Item(115).Hidden = $true
$ws.Rows.Item(117).Hidden = $true
$ws.Rows.Item(131).Hidden = $true
$ws.Rows.Item(136).Hidden = $true

# Re-apply the AutoFilter so its range grows to the new used range (the
# sheet dimension already reports A1:C136) and the _FilterDatabase defined
# name follows along.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:C136").AutoFilter(2, @("no comenzado"), 7)

# Row 101 ("en proceso") stays visible even though it no longer matches the
# "no comenzado" filter criteria, same as row 131 was before the edit -
# reapplying the filter recalculates hidden rows, so restore it explicitly.
$ws.Rows.Item(101).Hidden = $false

$fdb = $wb.Names.Item("Hoja1!_FilterDatabase")
$fdb.RefersTo = "=Hoja1!`$A`$1:`$C`$136"

# Move the active selection down one row to reflect the newly-added last row.
[void]$ws.Range("B139").Select()
